$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E2" = 3
    "F2" = 1
    "G2" = 14.597496
    "H2" = 43.792488
    "I2" = 0.3001995538392804
    "J2" = 0.3001995538392803
    "K2" = 3
    "L2" = 1
    "M2" = 14.597496
    "N2" = 43.792488
    "O2" = 0.3001995538392804
    "P2" = 0.3001995538392803
    "Q2" = 213.086889470016
    "R2" = 1917.782005230144
    "S2" = 0.09011977212530299
    "T2" = 0.09011977212530296
    "E3" = 3
    "F3" = 1
    "G3" = 14.597496
    "H3" = 43.792488
    "I3" = 0.3001995538392804
    "J3" = 0.3001995538392803
    "O3" = 0.09548910844461217
    "P3" = 0.09548910844461216
    "Q3" = 67.779837899496
    "R3" = 610.018541095464
    "S3" = 0.02866578775158323
    "T3" = 0.02866578775158322
    "E4" = 3
    "F4" = 1
    "G4" = 14.597496
    "H4" = 43.792488
    "I4" = 0.3001995538392804
    "J4" = 0.3001995538392803
    "M4" = 29.385228
    "N4" = 88.15568399999999
    "O4" = 0.6043113377161075
    "P4" = 0.6043113377161075
    "Q4" = 428.950748189088
    "R4" = 3860.556733701791
    "S4" = 0.1814139939623941
    "T4" = 0.1814139939623941
    "I5" = 0.09548910844461217
    "J5" = 0.09548910844461216
    "K5" = 3
    "L5" = 1
    "M5" = 14.597496
    "N5" = 43.792488
    "O5" = 0.3001995538392804
    "P5" = 0.3001995538392803
    "Q5" = 67.779837899496
    "R5" = 610.018541095464
    "S5" = 0.02866578775158323
    "T5" = 0.02866578775158322
    "I6" = 0.09548910844461217
    "J6" = 0.09548910844461216
    "O6" = 0.09548910844461217
    "P6" = 0.09548910844461216
    "S6" = 0.009118169831546903
    "T6" = 0.0091181698315469
    "I7" = 0.09548910844461217
    "J7" = 0.09548910844461216
    "M7" = 29.385228
    "N7" = 88.15568399999999
    "O7" = 0.6043113377161075
    "P7" = 0.6043113377161075
    "Q7" = 136.442989296228
    "R7" = 1227.986903666052
    "S7" = 0.05770515086148204
    "T7" = 0.05770515086148202
    "G8" = 29.385228
    "H8" = 88.15568399999999
    "I8" = 0.6043113377161075
    "J8" = 0.6043113377161075
    "K8" = 3
    "L8" = 1
    "M8" = 14.597496
    "N8" = 43.792488
    "O8" = 0.3001995538392804
    "P8" = 0.3001995538392803
    "Q8" = 428.950748189088
    "R8" = 3860.556733701791
    "S8" = 0.1814139939623941
    "T8" = 0.1814139939623941
    "G9" = 29.385228
    "H9" = 88.15568399999999
    "I9" = 0.6043113377161075
    "J9" = 0.6043113377161075
    "O9" = 0.09548910844461217
    "P9" = 0.09548910844461216
    "Q9" = 136.442989296228
    "R9" = 1227.986903666052
    "S9" = 0.05770515086148204
    "T9" = 0.05770515086148202
    "G10" = 29.385228
    "H10" = 88.15568399999999
    "I10" = 0.6043113377161075
    "J10" = 0.6043113377161075
    "M10" = 29.385228
    "N10" = 88.15568399999999
    "O10" = 0.6043113377161075
    "P10" = 0.6043113377161075
    "Q10" = 863.4916246119839
    "R10" = 7771.424621507855
    "S10" = 0.3651921928922313
    "T10" = 0.3651921928922313
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
